# Remove the Creative Commons Attribution-NonCommercial license text and
# the accompanying INCLUDEPICTURE field/logo image from the footer area of
# the document, leaving a single, empty paragraph behind (the template is
# now "clean" and no longer programmatically stamps a CC license).

$d = $word.ActiveDocument

# Locate the paragraph that carries the "Creative Commons ..." sentence and
# the paragraph immediately after it, which holds the INCLUDEPICTURE field
# (begin/instrText/separate/drawing/end) that renders the CC logo.
$ccTextParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Creative Commons*") {
        $ccTextParaIndex = $i
    }
}

if ($ccTextParaIndex -gt 0) {

    # The very next paragraph holds the INCLUDEPICTURE field + drawing.
    $ccFieldParaIndex = $ccTextParaIndex + 1
    if ($ccFieldParaIndex -le $d.Paragraphs.Count) {
        $pField = $d.Paragraphs($ccFieldParaIndex)
        # Delete the whole paragraph, including its own paragraph mark, so
        # the field codes / INCLUDEPICTURE drawing and the paragraph itself
        # disappear completely.
        $pField.Range.Delete()
    }

    # Now strip the license sentence's text runs, but keep the (now empty)
    # paragraph and its formatting (Footer style, centered, 18pt) intact.
    $pText = $d.Paragraphs($ccTextParaIndex)
    $rngText = $d.Range($pText.Range.Start, $pText.Range.End - 1)
    $rngText.Delete()
}
